$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 318, pushing the existing rows 318..375
# down to 319..376 (dimension grows from A1:R375 to A1:R376).
$ws.Rows("318").Insert()

# Populate the newly inserted row 318 with the new daily price record.
$ws.Cells.Item(318, 1).Value = 3
$ws.Cells.Item(318, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(318, 3).Value = "Coquimbo"
$ws.Cells.Item(318, 4).Value = 44476
$ws.Cells.Item(318, 5).Value = 5
$ws.Cells.Item(318, 6).Value = 100112024
$ws.Cells.Item(318, 7).Value = "Choclo"
$ws.Cells.Item(318, 8).Value = "Dulce o Americano"
$ws.Cells.Item(318, 9).Value = "Primera"
$ws.Cells.Item(318, 10).Value = 35
$ws.Cells.Item(318, 11).Value = 40000
$ws.Cells.Item(318, 12).Value = 40000
$ws.Cells.Item(318, 13).Value = 40000
$ws.Cells.Item(318, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(318, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(318, 16).Value = 571
$ws.Cells.Item(318, 17).Value = 70
$ws.Cells.Item(318, 18).Value = "Hortaliza"
